# "excel to json file"
# Insert a new first column ("country") in front of the existing data,
# shifting the current A/B/C columns to B/C/D, and fill the new column
# with the per-row country codes. Finally move the selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 3 columns of data one slot to the right (work from the
# rightmost column back to the left so we never clobber a value before it
# has been copied).
for ($r = 1; $r -le 3; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 1).Value()
}

# Populate the new leading column with the country data.
$ws.Range("A2").Value = "bdo_kr"
$ws.Range("A3").Value = "bdo_en"
$ws.Range("A1").Value = "country"

# Update the active selection to A2.
[void]$ws.Range("A2").Select()
